$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 81, shifting existing rows 81-111 down to 82-112.
$ws.Rows.Item(81).Insert()

# Populate the newly inserted row 81 with the new record.
$ws.Cells.Item(81, 1).Value = 5
$ws.Cells.Item(81, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(81, 3).Value = "Maule"
$ws.Cells.Item(81, 4).Value = 45205
$ws.Cells.Item(81, 5).Value = 7
$ws.Cells.Item(81, 6).Value = 300000000
$ws.Cells.Item(81, 7).Value = "Espárragos"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 4000
$ws.Cells.Item(81, 11).Value = 1100
$ws.Cells.Item(81, 12).Value = 1200
$ws.Cells.Item(81, 13).Value = 1150
$ws.Cells.Item(81, 14).Value = '$/kilo'
$ws.Cells.Item(81, 15).Value = "Región del Maule"
$ws.Cells.Item(81, 16).Value = 1150
$ws.Cells.Item(81, 17).Value = 1
$ws.Cells.Item(81, 18).Value = "Hortaliza"
